$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = @'
<Bold>e001 Fourth Armor Division Campaign</Bold> <InlineUIContainer><Button Content='r1.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
The campaign game of <Bold>Patton' Best</Bold> recreates the actions of the 4th Armored Division from late July 1944 through April 1945. 
<LineBreak/><LineBreak/>
Each day begins with a check of the Combat <InlineUIContainer><Button Content='Calendar' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> to see
 what the Division was doing on that day. The four possibilities are Refitting <InlineUIContainer><Button Content='r27.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>, 
an Advance scenario <InlineUIContainer><Button Content='r20.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>, 
a Battle scenario <InlineUIContainer><Button Content='r20.3' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>, 
or a Counterattack scenario <InlineUIContainer><Button Content='r20.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.  Click image to continue.
<LineBreak/><LineBreak/>
                                            <InlineUIContainer><Image Name='Nothing' Source='../images/Nothing.gif' Height='100' Width='100'></Image></InlineUIContainer>
'@

$ws.Range("B2").Value = $newText

$ws.Columns.Item(1).ColumnWidth = 7.35
$ws.Columns.Item(2).ColumnWidth = 180.5

$ws.Rows.Item(1).RowHeight = 165
$ws.Rows.Item(2).RowHeight = 165
$ws.Rows.Item(3).RowHeight = 28.5
$ws.Rows.Item(4).RowHeight = 28.5
$ws.Rows.Item(5).RowHeight = 128.45
$ws.Rows.Item(6).RowHeight = 60
$ws.Rows.Item(7).RowHeight = 90

$ws.Range("D2").Select()
